$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-04 Tuesday" "2025-03-05 Wednesday"

Replace-Text "40×53=" "19×16="
Replace-Text "30×32=" "54×77="
Replace-Text "48×17=" "40×77="
Replace-Text "59×81=" "49×25="
Replace-Text "63×62=" "30×29="

Replace-Text "29×83=" "44×11="
Replace-Text "12×75=" "42×83="
Replace-Text "12×17=" "18×52="
Replace-Text "16×16=" "56×49="
Replace-Text "46×38=" "16×19="

Replace-Text "94×25=" "87×43="
Replace-Text "89×11=" "32×39="
Replace-Text "92×79=" "76×20="
Replace-Text "21×21=" "60×69="
Replace-Text "28×41=" "41×45="

Replace-Text "82×93=" "98×68="
Replace-Text "71×69=" "22×88="
Replace-Text "55×33=" "55×16="
Replace-Text "61×30=" "46×69="
Replace-Text "88×87=" "49×95="

Replace-Text "69×62=" "63×19="
Replace-Text "80×45=" "59×20="
Replace-Text "27×12=" "84×86="
Replace-Text "25×49=" "89×99="
Replace-Text "19×32=" "87×22="
